$d = $word.ActiveDocument

# Edit 2: "A arrow on the left and right indicates" -> "A slider indicates"
$d.Content.Find.Execute("arrow on the left and right indicates", $true, $false, $false, $false, $false,
                         $true, 1, $false, "slider indicates", 2)
